# Equalize the population bars/gridlines that feed the "global estimations"
# chart so their scale matches the "localities estimations" chart. The chart
# is built from plain rectangles/lines grouped into a single top-level group
# shape on the slide; each bar/gridline/tick-label keeps its own name
# (rcNN / plNN / txNN) inside that group, so we reach them via GroupItems.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

$sh = $g.GroupItems.Item("rc11")
$sh.Left = 311.53826904296875
$sh.Top = 28.855039596557617
$sh.Width = 0.8488976955413818
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc12")
$sh.Left = 292.8624572753906
$sh.Top = 51.328269958496094
$sh.Width = 19.52472496032715
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc13")
$sh.Left = 244.8145751953125
$sh.Top = 73.80149841308594
$sh.Width = 67.57252502441406
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc14")
$sh.Left = 152.62387084960938
$sh.Top = 96.2748031616211
$sh.Width = 159.76332092285156
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc15")
$sh.Left = 148.5491485595703
$sh.Top = 118.74803161621094
$sh.Width = 163.83804321289062
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc23")
$sh.Width = 117.82756042480469
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc24")
$sh.Width = 19.015356063842773
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc25")
$sh.Left = 296.5975646972656
$sh.Top = 216.14654541015625
$sh.Width = 15.789527893066406
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc26")
$sh.Left = 279.1102600097656
$sh.Top = 238.61976623535156
$sh.Width = 33.27693176269531
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("rc27")
$sh.Left = 274.1866149902344
$sh.Top = 261.093017578125
$sh.Width = 38.20055389404297
$sh.Height = 20.225906372070312

$sh = $g.GroupItems.Item("pl32")
$sh.Left = 142.60678100585938
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl33")
$sh.Left = 185.05189514160156
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl34")
$sh.Left = 227.49700927734375
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl35")
$sh.Left = 269.94207763671875
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl37")
$sh.Left = 354.8323059082031
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl38")
$sh.Left = 397.27734375
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl39")
$sh.Left = 439.7224426269531
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("pl40")
$sh.Left = 482.16748046875
$sh.Top = 284.68994140625

$sh = $g.GroupItems.Item("tx41")
$sh.Left = 129.27560424804688
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx42")
$sh.Left = 171.72064208984375
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx43")
$sh.Left = 214.16575622558594
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx44")
$sh.Left = 256.6108093261719
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx46")
$sh.Left = 343.6322937011719
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx47")
$sh.Left = 386.07733154296875
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx48")
$sh.Left = 428.5224609375
$sh.Top = 289.27142333984375

$sh = $g.GroupItems.Item("tx49")
$sh.Left = 470.9674987792969
$sh.Top = 289.27142333984375
